# Add "coming soon" placeholder message to the Charts tab and make it the
# active sheet (it was previously empty with the Metadata tab selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Charts")

$ws.Range("A1").Value = "Automatically generated chart(s) coming soon to this tab."

# Make "Charts" the selected/active sheet (was "Metadata").
$ws.Activate()
